$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, derived from the day-over-day market refresh
$updates = @{
    'D2' = '28.445.18'
    'E2' = '  +3.50%  '
    'D3' = '1.591.29'
    'E3' = '  +1.54%  '
    'E4' = '  +0.90%  '
    'D5' = '213.08'
    'E5' = '  +0.89%  '
    'E6' = '  +0.35%  '
    'E7' = '  +0.90%  '
    'D8' = '24.44'
    'E8' = '  +7.90%  '
    'E9' = '  +0.44%  '
    'E10' = '  +0.89%  '
    'D11' = '0.0886'
    'E11' = '  +1.68%  '
    'D12' = '1.817.10'
    'E12' = '  +1.51%  '
    'D13' = '1.602.07'
    'E13' = '  +2.34%  '
    'E14' = '  +2.20%  '
    'E15' = '  -0.05%  '
    'D16' = '28.461.63'
    'E16' = '  +3.66%  '
    'D17' = '63.10'
    'E17' = '  +1.17%  '
    'D18' = '229.75'
    'E18' = '  +1.93%  '
    'D19' = '0.0₃0707'
    'E19' = '  +0.49%  '
    'D20' = '7.47'
    'E20' = '  -0.19%  '
    'E21' = '  +0.84%  '
    'D22' = '4.06'
    'E22' = '  -1.17%  '
    'D23' = '9.34'
    'E23' = '  -0.44%  '
    'E24' = '  +0.85%  '
    'D25' = '151.67'
    'E25' = '  +1.09%  '
    'D26' = '15.22'
    'E26' = '  +0.51%  '
    'E27' = '  -0.61%  '
    'E28' = '  -0.90%  '
    'E29' = '  +0.94%  '
    'E30' = '  -0.59%  '
    'E31' = '  +0.33%  '
    'E32' = '  +0.55%  '
    'E33' = '  +0.79%  '
    'D34' = '1.401.25'
    'E34' = '  -3.30%  '
    'E35' = '  -0.53%  '
    'E36' = '  -9.06%  '
    'E37' = '  +0.85%  '
    'B38' = 'MXToken'
    'C38' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D38' = '2.58'
    'E38' = '  +8.92%  '
    'B39' = 'VeChain'
    'C39' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D39' = '0.0167'
    'E39' = '  -0.41%  '
    'E40' = '  +0.37%  '
    'D41' = '0.812'
    'E41' = '  -0.11%  '
    'E42' = '  +0.86%  '
    'E43' = '  -2.15%  '
    'E44' = '  +0.20%  '
    'D45' = '0.981'
    'E45' = '  +0.65%  '
    'E46' = '  -1.45%  '
    'D47' = '1.726.56'
    'E47' = '  +1.44%  '
    'E48' = '  +1.77%  '
    'D49' = '87.23'
    'E49' = '  +0.42%  '
    'E50' = '  +0.62%  '
    'E51' = '  -1.01%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (e.g. "213.08", "63.10") keep their exact formatting
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
